# Append/update run: 2025-10-06 06:27 JST
# Updates the "ランサーズ" (case list) sheet: refreshes the timestamp on the
# still-present rows, swaps in the newest 4 listings (rows 2,4,5 get new
# job postings; row 3 is untouched apart from its timestamp), and drops the
# previously-captured rows 6-9 which fell out of the refreshed window.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newStamp = "2025-10-06 06:27:24"

# --- Row 2: new top listing -------------------------------------------------
$ws.Range("A2").Value = $newStamp
$ws.Range("B2").Value = "あなたAIクローン構築パートナー募集・モデル制作&新規依頼"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5399534"
$ws.Range("G2").Value = 303
$ws.Range("H2").Value = "🔥AI,Ai"

# --- Row 3: timestamp refresh only ------------------------------------------
$ws.Range("A3").Value = $newStamp

# --- Row 4: new listing, and it loses its skill-summary cell ---------------
$ws.Range("A4").Value = $newStamp
$ws.Range("B4").Value = "【急募】AWSマイクロサービスのデバッグ・最適化支援者募集"
$ws.Range("D4").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5407390"
$ws.Range("G4").Value = 10
$ws.Range("H4").ClearContents()

# --- Row 5: new listing, and it loses its skill-summary cell ---------------
$ws.Range("A5").Value = $newStamp
$ws.Range("B5").Value = "【急募】サーバー移転後のWelcartクレカ決済不具合解消依頼"
$ws.Range("D5").Value = "~ 5,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5407516"
$ws.Range("G5").Value = 10
$ws.Range("H5").ClearContents()

# --- Drop rows 6-9 (no longer part of the captured window) -----------------
$ws.Range("A6:H9").EntireRow.Delete()

# --- Hyperlinks: this engine's Hyperlinks.Delete() clears the whole sheet's
# collection regardless of the range it's called on, so rebuild the four
# that should remain (F2:F5) with their (partly updated) targets after
# wiping the stale set (which still pointed rows 2/4/5 at their old URLs
# and rows 6-9 at URLs that no longer exist on the sheet).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5399534")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5407281")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5407390")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5407516")

# Hyperlinks.Add() re-stamps the font but minted a fresh cell style rather
# than reusing the workbook's existing "Hyperlink" cell style; point the
# cells back at the named style so F2:F5 keep looking like the original
# hyperlink-styled cells.
$ws.Range("F2:F5").Style = "Hyperlink"

# --- Column widths -----------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 42
$ws.Columns.Item(4).ColumnWidth = 28
